$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a single date value (serial 46061) that was
# updated by the automated export to the following day (serial 46062) for
# every data row (rows 2 through 515).
$ws.Range("C2:C515").Value = 46062
